$wb = $excel.ActiveWorkbook

# =========================================================================
# Summary sheet
# =========================================================================
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 1300
$summary.Range("B3").Value = 1300.23
$summary.Range("B4").Value = 0.23
$summary.Range("B5").Value = 0.42
$summary.Range("B6").Value = 11
$summary.Range("B7").Value = 6
$summary.Range("B9").Value = 54.55
$summary.Range("B11").Value = 13

# =========================================================================
# Strategy Status sheet - EMAArbitrage strategy row removed (row 3), the
# rest of the table shifts up, and the MarketMaking row (now row 4) gets
# its stats refreshed with the newest trading results
# =========================================================================
$status = $wb.Worksheets.Item("Strategy Status")
$status.Rows.Item(3).Delete()
$status.Range("C4").Value = 100.23
$status.Range("D4").Value = 11
$status.Range("E4").Value = 0.23
$status.Range("F4").Value = 0.23
$status.Range("G4").Value = 54.55

# =========================================================================
# Trade-log sheets ("All Trades" and "MarketMaking" mirror each other).
# Only touch the cells that actually changed on the six existing rows,
# then append five brand-new closed trades as rows 8-12.
# =========================================================================
$tradeSheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $tradeSheetNames) {
    $ts = $wb.Worksheets.Item($sheetName)

    # --- Row 2 (Trade #1) ---
    $ts.Range("C2").Value = "19:43:29"
    $ts.Range("E2").Value = "UP"
    $ts.Range("F2").Value = 0.41
    $ts.Range("G2").Value = 0.4
    $ts.Range("I2").Value = -2.439
    $ts.Range("J2").Value = -0.01
    $ts.Range("K2").Value = 99.98999999999999
    $ts.Range("Q2").Value = 0.11

    # --- Row 3 (Trade #2) ---
    $ts.Range("C3").Value = "19:43:37"
    $ts.Range("F3").Value = 0.4
    $ts.Range("G3").Value = 0.37
    $ts.Range("I3").Value = -7.5
    $ts.Range("J3").Value = -0.03
    $ts.Range("K3").Value = 99.95999999999999

    # --- Row 4 (Trade #3) ---
    $ts.Range("C4").Value = "19:43:45"
    $ts.Range("F4").Value = 0.64
    $ts.Range("G4").Value = 0.7
    $ts.Range("I4").Value = 9.375
    $ts.Range("J4").Value = 0.06
    $ts.Range("K4").Value = 100.02
    $ts.Range("Q4").Value = 0.14

    # --- Row 5 (Trade #4) ---
    $ts.Range("C5").Value = "19:43:53"
    $ts.Range("E5").Value = "DOWN"
    $ts.Range("F5").Value = 0.7
    $ts.Range("G5").Value = 0.72
    $ts.Range("I5").Value = 2.8571
    $ts.Range("J5").Value = 0.02
    $ts.Range("K5").Value = 100.04
    $ts.Range("Q5").Value = 0.13

    # --- Row 6 (Trade #5) ---
    $ts.Range("C6").Value = "19:44:01"
    $ts.Range("E6").Value = "DOWN"
    $ts.Range("F6").Value = 0.7272729999999999
    $ts.Range("G6").Value = 0.78
    $ts.Range("I6").Value = 7.25
    $ts.Range("J6").Value = 0.05
    $ts.Range("K6").Value = 100.09

    # --- Row 7 (Trade #6) : was the still-OPEN trade, now CLOSED ---
    $ts.Range("C7").Value = "19:44:10"
    $ts.Range("E7").Value = "UP"
    $ts.Range("F7").Value = 0.22
    $ts.Range("G7").Value = 0.1
    $ts.Range("H7").Value = "CLOSED"
    $ts.Range("I7").Value = -54.5455
    $ts.Range("J7").Value = -0.12
    $ts.Range("K7").Value = 99.97
    $ts.Range("P7").Value = "early_exit"
    $ts.Range("Q7").Value = 0.14

    # --- New rows 8-12 : five freshly-appended closed trades ---
    $newRows = @(
        @(7,  "19:44:18", "DOWN", 0.93, 0.95,               2.1505,   0.02,  99.98999999999999),
        @(8,  "19:44:33", "UP",   0.03, 0.03,                0,        0,     99.98999999999999),
        @(9,  "19:44:40", "DOWN", 0.97, 0.98,                1.0309,   0.01,  100),
        @(10, "19:47:27", "UP",   0.83, 0.8100000000000001, -2.4096,  -0.02,  99.98),
        @(11, "19:47:35", "DOWN", 0.19, 0.44,                131.5789, 0.25,  100.23)
    )
    $qVals = @(0.15, 0.11, 0.13, 0.11, 0.14)

    $r = 8
    for ($i = 0; $i -lt $newRows.Count; $i++) {
        $row = $newRows[$i]
        $ts.Cells.Item($r, 1).Value  = $row[0]                      # Trade #

        $ts.Cells.Item($r, 2).NumberFormat = "@"
        $ts.Cells.Item($r, 2).Value  = "2026-02-17"                 # Date

        $ts.Cells.Item($r, 3).Value  = $row[1]                      # Time
        $ts.Cells.Item($r, 4).Value  = "MarketMaking"                # Strategy
        $ts.Cells.Item($r, 5).Value  = $row[2]                      # Side
        $ts.Cells.Item($r, 6).Value  = $row[3]                      # Entry Price
        $ts.Cells.Item($r, 7).Value  = $row[4]                      # Exit Price
        $ts.Cells.Item($r, 8).Value  = "CLOSED"                      # Status
        $ts.Cells.Item($r, 9).Value  = $row[5]                      # P&L %
        $ts.Cells.Item($r, 10).Value = $row[6]                      # P&L $
        $ts.Cells.Item($r, 11).Value = $row[7]                      # Capital After
        $ts.Cells.Item($r, 12).Value = 0                             # Entry Slippage
        $ts.Cells.Item($r, 13).Value = 0                             # Exit Slippage
        $ts.Cells.Item($r, 14).Value = 0.6                           # Confidence
        $ts.Cells.Item($r, 15).Value = "Normal spread capture: 19600 bps"  # Entry Reason
        $ts.Cells.Item($r, 16).Value = "early_exit"                  # Exit Reason
        $ts.Cells.Item($r, 17).Value = $qVals[$i]                    # Duration (min)

        $r = $r + 1
    }
}
